$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (佐藤茂邸 / 外壁 / 外壁 / 大工 / 板張り / 000 / 0 / ー),
# shifting every subsequent row up by one.
$ws.Rows.Item(2).Delete()

# Remove the row that is now "松田邸" (previously row 10), shifting the
# remaining rows up by one again.
$ws.Rows.Item(9).Delete()

# The row that used to be "佐藤茂邸 / 内壁 ..." (old row 11) is now row 9.
# Overwrite it with the new finish data.
$ws.Range("A9").Value = "佐藤茂邸"
$ws.Range("B9").Value = "外壁"
$ws.Range("C9").Value = "外壁"
$ws.Range("D9").Value = "アイカ工業"
$ws.Range("E9").Value = "ジョリパッドネオ∞ JQ-620"
$ws.Range("F9").Value = "ー"
$ws.Range("G9").Value = "T5005"
$ws.Range("H9").Value = ""
